# Update Name of Algo
# Apply updated numeric results to Sheet1 (result_data_RandomForest)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.259999999999996
$ws.Range("B3").Value = 6.052400000000005
$ws.Range("B14").Value = 5.900900000000002
$ws.Range("B21").Value = 9.423600000000006
$ws.Range("B23").Value = 8.9777
$ws.Range("B25").Value = 5.535400000000002
$ws.Range("D25").Value = -8.201200000000002
$ws.Range("B26").Value = 5.626200000000003
$ws.Range("D27").Value = -8.868600000000004
$ws.Range("B29").Value = 5.009200000000003
$ws.Range("D31").Value = -8.535800000000004
$ws.Range("D39").Value = -8.004399999999999
$ws.Range("D48").Value = -7.478699999999995
$ws.Range("D51").Value = -7.796199999999998
$ws.Range("D52").Value = -7.792999999999999
$ws.Range("B53").Value = 5.355099999999998
$ws.Range("D55").Value = -8.359999999999998
$ws.Range("D56").Value = -7.891299999999996
$ws.Range("B57").Value = 5.068799999999995
$ws.Range("D57").Value = -8.182499999999999
$ws.Range("B59").Value = 4.959499999999998
$ws.Range("B69").Value = 5.276499999999995
$ws.Range("D73").Value = -7.934899999999996
$ws.Range("B79").Value = 9.493500000000006
$ws.Range("B83").Value = 5.308099999999997
$ws.Range("D89").Value = -5.969299999999999
$ws.Range("D90").Value = -8.108400000000003
$ws.Range("B91").Value = 4.908199999999999
$ws.Range("D92").Value = -6.007099999999999
$ws.Range("B93").Value = 5.880700000000005
